# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handback DateTime"
# timestamps on the handback status report to reflect the new report run.

$wb = $excel.ActiveWorkbook

# Overview sheet: Latest HO Xliff Generate Date for the 41adc32b... file
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-09-02 00:55:43"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-09-02 00:55:39"
$wsZhCn.Range("K4").Value = "2016-09-02 00:56:00"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-09-02 00:55:43"
$wsDeDe.Range("K4").Value = "2016-09-02 00:56:16"
